# Update the "segment-pivot" formatted pivot table with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("segment-pivot")

# Update the data values for rows 3-5 (columns C:H) to reflect the refreshed pivot numbers.
$ws.Range("C3").Value = 779
$ws.Range("D3").Value = 426
$ws.Range("E3").Value = 356
$ws.Range("F3").Value = 339
$ws.Range("G3").Value = 390
$ws.Range("H3").Value = 12136

$ws.Range("C4").Value = 680
$ws.Range("D4").Value = 220
$ws.Range("E4").Value = 165
$ws.Range("F4").Value = 165
$ws.Range("G4").Value = 168
$ws.Range("H4").Value = 13028

$ws.Range("C5").Value = 507
$ws.Range("D5").Value = 120
$ws.Range("E5").Value = 77
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 80
$ws.Range("H5").Value = 13558

# Move the active cell selection to B3, matching the saved view state.
$ws.Range("B3").Select()
